$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 32 for New York state hospitalization data, 15 April 2020 (serial 43936)
# Copy the date format from the cell above first, then set the value.
$ws.Cells.Item(31, 1).Copy()
$ws.Cells.Item(32, 1).PasteSpecial(-4122)
$ws.Cells.Item(32, 1).Value = 43936

$ws.Cells.Item(32, 2).Value = -606
$ws.Cells.Item(32, 3).Value = -154
$ws.Cells.Item(32, 4).Value = -62
$ws.Cells.Item(32, 6).Value = 606
$ws.Cells.Item(32, 7).Value = 1996

# Update selection to match post-edit state (F33)
$ws.Range("F33").Select()
